# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" worksheet (per-fund holdings detail, same
# layout as the existing 2022-Q2 / 2022-Q1 sheets) right after the
# "总计" (totals) sheet, and records the new quarter as an extra row
# at the top of the "总计" summary table.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Update the "总计" summary sheet: insert a new row right under the
#    header for 2022-Q3, pushing the existing 2022-Q2 / 2022-Q1 rows
#    down by one, then fill in all three data rows with their final
#    values.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Push old row 2 (2022-Q2) / row 3 (2022-Q1) down to rows 3 / 4.
$total.Rows.Item(3).Insert()

# Row 2: new 2022-Q3 entry.
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.11

# Row 3: 2022-Q2 entry (was row 2).
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 3
$total.Range("D3").Value = 0.29

# Row 4: 2022-Q1 entry (was row 3), re-index its leading counter.
$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.09

# Give the newly-inserted index cell (A3) the same bold / bordered /
# centered style used by the other index-column cells (A2, A4).
$idxCells = $total.Range("A2:A4")
$idxCells.Font.Bold = $true
$idxCells.Borders.LineStyle = 1
$idxCells.HorizontalAlignment = -4108
$idxCells.VerticalAlignment = -4160

# ------------------------------------------------------------------
# 2) Insert a brand-new worksheet "2022-Q3" right after "总计" holding
#    the per-fund holdings detail, mirroring the layout used by the
#    2022-Q2 / 2022-Q1 sheets.
# ------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$q3.Range("A2").Value = 0
$q3.Range("B2").NumberFormat = "@"
$q3.Range("B2").Value = "014275"
$q3.Range("C2").Value = "易方达北交所精选两年定开混合A"
$q3.Range("D2").NumberFormat = "@"
$q3.Range("D2").Value = "3.58"
$q3.Range("E2").NumberFormat = "@"
$q3.Range("E2").Value = "61.75"
$q3.Range("F2").NumberFormat = "@"
$q3.Range("F2").Value = "2.34"
$q3.Range("G2").NumberFormat = "@"
$q3.Range("G2").Value = "0.0838"
$q3.Range("H2").Value = 10

$q3.Range("A3").Value = 1
$q3.Range("B3").NumberFormat = "@"
$q3.Range("B3").Value = "014276"
$q3.Range("C3").Value = "易方达北交所精选两年定开混合C"
$q3.Range("D3").NumberFormat = "@"
$q3.Range("D3").Value = "0.92"
$q3.Range("E3").NumberFormat = "@"
$q3.Range("E3").Value = "61.75"
$q3.Range("F3").NumberFormat = "@"
$q3.Range("F3").Value = "2.34"
$q3.Range("G3").NumberFormat = "@"
$q3.Range("G3").Value = "0.0215"
$q3.Range("H3").Value = 10

# Match the bold/bordered header + index-column styling used on the
# other detail sheets (2022-Q2 / 2022-Q1).
$headerCells = $q3.Range("B1:H1")
$headerCells.Font.Bold = $true
$headerCells.Borders.LineStyle = 1
$headerCells.HorizontalAlignment = -4108
$headerCells.VerticalAlignment = -4160

$q3IdxCells = $q3.Range("A2:A3")
$q3IdxCells.Font.Bold = $true
$q3IdxCells.Borders.LineStyle = 1
$q3IdxCells.HorizontalAlignment = -4108
$q3IdxCells.VerticalAlignment = -4160

# The original "2022-Q1" sheet was the last / active tab before this
# edit; keep that same sheet active afterwards.
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Activate()

Write-Output "done"
